$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Survey 3" data row (spreadsheet row 4) that was previously missing
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 18

# Match the author's final selection on the newly-entered row
$ws.Range("B4:D4").Select()
